$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("在D盘根目录下放置", $true, $false, $false, $false, $false, $true, 1, $false, "在项目目录\src\main\resources下放置", 2)
